$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Final desired data for A2:C17 (Miles McBride / New York Knicks row removed,
# remaining players reordered and Jaden Ivey now on Detroit Pistons).
$data = @(
    @("Jaden Ivey", "PG,SG", "Detroit Pistons"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Aaron Gordon", "PF,C", "Denver Nuggets"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Khris Middleton", "SF", "Milwaukee Bucks"),
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove the now-obsolete last row (previously row 18) so the table is 17 rows
# (1 header + 16 data rows) instead of 18.
$ws.Rows.Item(18).Delete()
